$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B5 to use the new Arabic string for "male"
$ws.Range("B5").Value = "الذكر"

# Autofit column B to match the bestFit width recorded after editing
$ws.Columns.Item(2).AutoFit() | Out-Null

# Select D16 to match the saved cursor position
$ws.Range("D16").Select() | Out-Null

# Configure page setup as recorded in the saved workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
